# Remove mock data from the "取得" (Get) worksheet's asset table.
# The table (Table14, A1:F101) keeps its headers and structure, but every
# sample/mock row (2-42) loses its data:
#   - Columns C (アセットID), D (アセット名), E (タイプ) are fully cleared,
#     including their cell formatting, so the cells disappear entirely.
#   - Columns A (組織単位ID), B (組織単位名), F (値) keep their existing
#     cell style but are left blank.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("取得")

# Clear content + formatting for the columns that should vanish entirely.
$ws.Range("D2:E42").Clear()

# Blank out all six columns of mock data rows (keeps formatting on A/B/F).
$ws.Range("A2:F42").ClearContents()

# A couple of rows had a taller custom height because of the long mock
# text that used to live in column D/E; restore the default auto height
# now that the content is gone.
$ws.Range("A2:F42").EntireRow.AutoFit()
